$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "is_closed"

# "true"/"false" look like booleans to Excel's Value setter, so build them
# as text formulas first, then convert each to a plain value in place -
# this keeps the cell as a shared-string text cell (t="s") without
# picking up a quote-prefix style.
$ws.Range("B2").Value = "s"
$ws.Range("C2").Formula = "=""true"""
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)

$ws.Range("B3").Value = "склад"
$ws.Range("C3").Formula = "=""true"""
$ws.Range("C3").Copy()
$ws.Range("C3").PasteSpecial(-4163)

$ws.Range("B4").Value = "storage"
$ws.Range("C4").Formula = "=""false"""
$ws.Range("C4").Copy()
$ws.Range("C4").PasteSpecial(-4163)

$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(3).EntireColumn.AutoFit()
